$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the measurement values in row 5 to 2 decimal places (custom accuracy)
$ws.Range("B5").Value = 5.51
$ws.Range("C5").Value = 4.22
$ws.Range("D5").Value = 0.06
$ws.Range("E5").Value = 12.34
$ws.Range("F5").Value = 9.93
$ws.Range("G5").Value = 4.12
$ws.Range("H5").Value = 23.16
$ws.Range("I5").Value = 6.83
$ws.Range("J5").Value = 3.43
$ws.Range("K5").Value = 4.56
$ws.Range("L5").Value = 5.02
$ws.Range("M5").Value = 5.33
$ws.Range("N5").Value = 1.43
$ws.Range("O5").Value = 4.49
$ws.Range("P5").Value = 6.28
$ws.Range("Q5").Value = 3.87
$ws.Range("R5").Value = 0.06
$ws.Range("S5").Value = 0.32
$ws.Range("T5").Value = 61.39
$ws.Range("U5").Value = 12.88
$ws.Range("V5").Value = 4.15
$ws.Range("W5").Value = 8.46
$ws.Range("X5").Value = 4.71
$ws.Range("Y5").Value = 0.74
$ws.Range("Z5").Value = 10.86
$ws.Range("AA5").Value = 3.6
$ws.Range("AB5").Value = 3.62
$ws.Range("AC5").Value = 4.3
$ws.Range("AD5").Value = 5.49
$ws.Range("AE5").Value = 0.33
$ws.Range("AF5").Value = 21.33
$ws.Range("AG5").Value = 2.55
$ws.Range("AH5").Value = 5.04

# Remove the last data row (row 6) entirely
$ws.Rows.Item(6).Delete()
